# [Fonds de solidarite] Add 2020-12-16 data
# Updates nombre_aides (C), nombre_entreprises (D) and montant_total (E)
# for the rows whose section/region totals changed with the new data.
#
# The source cells are stored as text (inline strings), so we force the
# NumberFormat to "@" (Text) before writing the value - this keeps Excel
# from re-interpreting the numeric-looking strings as actual numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 21 - Bourgogne-Franche-Comte / F Construction
Set-TextValue "C21" "49"
Set-TextValue "E21" "141171.00"

# Row 24 - Bourgogne-Franche-Comte / I Hebergement et restauration
Set-TextValue "C24" "519"
Set-TextValue "D24" "430"
Set-TextValue "E24" "4439957.92"

# Row 28 - Bourgogne-Franche-Comte / M Activites specialisees, scientifiques et techniques
Set-TextValue "C28" "58"
Set-TextValue "E28" "350148.39"

# Row 32 - Bourgogne-Franche-Comte / R Arts, spectacles et activites recreatives
Set-TextValue "C32" "100"
Set-TextValue "E32" "1219245.89"

# Row 84 - Grand Est / I Hebergement et restauration
Set-TextValue "C84" "846"
Set-TextValue "E84" "7693397.41"

# Row 92 - Grand Est / R Arts, spectacles et activites recreatives
Set-TextValue "C92" "132"
Set-TextValue "E92" "1221799.93"

# Row 130 - Hauts-de-France / I Hebergement et restauration
Set-TextValue "C130" "1162"
Set-TextValue "E130" "9866904.14"

# Row 134 - Hauts-de-France / M Activites specialisees, scientifiques et techniques
Set-TextValue "C134" "170"
Set-TextValue "E134" "1695798.92"

# Row 147 - Ile-de-France / I Hebergement et restauration
Set-TextValue "C147" "5096"
Set-TextValue "E147" "35690046.99"

# Row 151 - Ile-de-France / M Activites specialisees, scientifiques et techniques
Set-TextValue "C151" "1642"
Set-TextValue "E151" "6317634.03"

# Row 155 - Ile-de-France / R Arts, spectacles et activites recreatives
Set-TextValue "C155" "861"
Set-TextValue "E155" "5245048.93"

# Row 163 - La Reunion / M Activites specialisees, scientifiques et techniques
Set-TextValue "C163" "16"
Set-TextValue "D163" "16"
Set-TextValue "E163" "170500.00"
